# Added work on progenitors to background and results sections.
# -> Append a new, blank slide at the end of the deck (slide 10),
#    using the same "Blank" layout used elsewhere in the deck.

$p = $ppt.ActivePresentation

# ppLayoutBlank = 12 -> maps to this deck's "Blank" custom layout
# (the same layout most of the existing slides are not on, but the
# one with no placeholders, matching a freshly-inserted empty slide).
$s = $p.Slides.Add($p.Slides.Count + 1, 12)
